$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.21
$ws.Cells.Item(2, 3).Value = 0.55
$ws.Cells.Item(2, 10).Value = 0.01
$ws.Cells.Item(2, 16).Value = 0.165
$ws.Cells.Item(2, 19).Value = 0.065

$ws.Cells.Item(3, 2).Value = 0.01724137931034483
$ws.Cells.Item(3, 3).Value = 0.0603448275862069
$ws.Cells.Item(3, 10).Value = 0.03448275862068965
$ws.Cells.Item(3, 16).Value = 0.646551724137931
$ws.Cells.Item(3, 19).Value = 0.2413793103448276

$ws.Cells.Item(4, 10).Value = 0.05714285714285714
$ws.Cells.Item(4, 16).Value = 0.7428571428571429
$ws.Cells.Item(4, 19).Value = 0.2

$ws.Cells.Item(5, 16).Value = 0.5
$ws.Cells.Item(5, 19).Value = 0.5

$ws.Cells.Item(6, 2).Value = 0.06111111111111111
$ws.Cells.Item(6, 4).Value = 0.02222222222222222
$ws.Cells.Item(6, 5).Value = 0.005555555555555556
$ws.Cells.Item(6, 6).Value = 0.08888888888888889
$ws.Cells.Item(6, 10).Value = 0.2333333333333333
$ws.Cells.Item(6, 15).Value = 0.02222222222222222
$ws.Cells.Item(6, 17).Value = 0.1666666666666667
$ws.Cells.Item(6, 18).Value = 0.07222222222222222
$ws.Cells.Item(6, 19).Value = 0.3277777777777778

$ws.Cells.Item(7, 2).Value = 0.08496732026143791
$ws.Cells.Item(7, 4).Value = 0.0196078431372549
$ws.Cells.Item(7, 5).Value = 0.006535947712418301
$ws.Cells.Item(7, 6).Value = 0.06535947712418301
$ws.Cells.Item(7, 10).Value = 0.0718954248366013
$ws.Cells.Item(7, 15).Value = 0.0130718954248366
$ws.Cells.Item(7, 17).Value = 0.1568627450980392
$ws.Cells.Item(7, 18).Value = 0.09803921568627451
$ws.Cells.Item(7, 19).Value = 0.4836601307189543

$ws.Cells.Item(8, 2).Value = 0.05298013245033113
$ws.Cells.Item(8, 4).Value = 0.01324503311258278
$ws.Cells.Item(8, 6).Value = 0.04856512141280353
$ws.Cells.Item(8, 10).Value = 0.108167770419426
$ws.Cells.Item(8, 15).Value = 0.01545253863134658
$ws.Cells.Item(8, 17).Value = 0.2097130242825607
$ws.Cells.Item(8, 18).Value = 0.08167770419426049
$ws.Cells.Item(8, 19).Value = 0.4701986754966888

$ws.Cells.Item(9, 2).Value = 0.04489795918367347
$ws.Cells.Item(9, 4).Value = 0.02040816326530612
$ws.Cells.Item(9, 6).Value = 0.05714285714285714
$ws.Cells.Item(9, 10).Value = 0.1142857142857143
$ws.Cells.Item(9, 15).Value = 0.0163265306122449
$ws.Cells.Item(9, 17).Value = 0.1306122448979592
$ws.Cells.Item(9, 18).Value = 0.08571428571428572
$ws.Cells.Item(9, 19).Value = 0.5306122448979592

$ws.Cells.Item(10, 2).Value = 0.07922535211267606
$ws.Cells.Item(10, 4).Value = 0.01496478873239437
$ws.Cells.Item(10, 6).Value = 0.05721830985915493
$ws.Cells.Item(10, 10).Value = 0.1258802816901408
$ws.Cells.Item(10, 15).Value = 0.01232394366197183
$ws.Cells.Item(10, 17).Value = 0.1892605633802817
$ws.Cells.Item(10, 18).Value = 0.1161971830985915
$ws.Cells.Item(10, 19).Value = 0.4049295774647887

$ws.Cells.Item(11, 7).Value = 0.0990990990990991
$ws.Cells.Item(11, 10).Value = 0.08108108108108109
$ws.Cells.Item(11, 11).Value = 0.1441441441441441
$ws.Cells.Item(11, 12).Value = 0.6621621621621622
$ws.Cells.Item(11, 19).Value = 0.01351351351351351

$ws.Cells.Item(12, 7).Value = 0.7666666666666667
$ws.Cells.Item(12, 10).Value = 0.1933333333333333
$ws.Cells.Item(12, 12).Value = 0.01333333333333333
$ws.Cells.Item(12, 19).Value = 0.02666666666666667

$ws.Cells.Item(13, 7).Value = 0.7058823529411765
$ws.Cells.Item(13, 10).Value = 0.2941176470588235

$ws.Cells.Item(15, 6).Value = 0.01052631578947368
$ws.Cells.Item(15, 8).Value = 0.1947368421052632
$ws.Cells.Item(15, 9).Value = 0.1052631578947368
$ws.Cells.Item(15, 10).Value = 0.3578947368421053
$ws.Cells.Item(15, 11).Value = 0.04736842105263158
$ws.Cells.Item(15, 13).Value = 0.005263157894736842
$ws.Cells.Item(15, 15).Value = 0.03684210526315789
$ws.Cells.Item(15, 19).Value = 0.2421052631578947

$ws.Cells.Item(16, 6).Value = 0.03053435114503817
$ws.Cells.Item(16, 8).Value = 0.1908396946564886
$ws.Cells.Item(16, 9).Value = 0.1297709923664122
$ws.Cells.Item(16, 10).Value = 0.3969465648854962
$ws.Cells.Item(16, 11).Value = 0.06870229007633588
$ws.Cells.Item(16, 13).Value = 0.007633587786259542
$ws.Cells.Item(16, 15).Value = 0.07633587786259542
$ws.Cells.Item(16, 19).Value = 0.09923664122137404

$ws.Cells.Item(17, 6).Value = 0.01256281407035176
$ws.Cells.Item(17, 8).Value = 0.1934673366834171
$ws.Cells.Item(17, 9).Value = 0.1080402010050251
$ws.Cells.Item(17, 10).Value = 0.4045226130653266
$ws.Cells.Item(17, 11).Value = 0.07788944723618091
$ws.Cells.Item(17, 13).Value = 0.01507537688442211
$ws.Cells.Item(17, 15).Value = 0.05276381909547739
$ws.Cells.Item(17, 19).Value = 0.135678391959799

$ws.Cells.Item(18, 6).Value = 0.009216589861751152
$ws.Cells.Item(18, 8).Value = 0.1336405529953917
$ws.Cells.Item(18, 9).Value = 0.119815668202765
$ws.Cells.Item(18, 10).Value = 0.4147465437788018
$ws.Cells.Item(18, 11).Value = 0.1152073732718894
$ws.Cells.Item(18, 13).Value = 0.009216589861751152
$ws.Cells.Item(18, 15).Value = 0.05069124423963134
$ws.Cells.Item(18, 19).Value = 0.1474654377880184

$ws.Cells.Item(19, 6).Value = 0.01472868217054264
$ws.Cells.Item(19, 8).Value = 0.2232558139534884
$ws.Cells.Item(19, 9).Value = 0.1085271317829457
$ws.Cells.Item(19, 10).Value = 0.3434108527131783
$ws.Cells.Item(19, 11).Value = 0.08682170542635659
$ws.Cells.Item(19, 13).Value = 0.01937984496124031
$ws.Cells.Item(19, 15).Value = 0.06666666666666667
$ws.Cells.Item(19, 19).Value = 0.1372093023255814
